$d = $word.ActiveDocument

# The paragraph "Esteu participant ... visibilitat del cel nocturn." is
# currently split across a dozen runs (each carrying the same rPr). Replace
# the whole paragraph's content with a single run with no direct character
# formatting, and tweak the wording around the constellation name at the
# same time ("constel·lació Perseus" -> " Constel·lació de Perseu").
foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    if ($r.Text -like "Esteu participant*visibilitat del cel nocturn.*") {
        # Exclude the trailing paragraph mark so only the run content is replaced.
        $rng = $d.Range($r.Start, $r.End - 1)

        $newText = "Esteu participant en una campanya mundial per observar i anotar la brillantor de les estrelles més febles que es poden veure, com a mitjà per mesurar la contaminació lumínica en un lloc determinat. Localitzant i observant la  Constel·lació de Perseu a la nit i comparant la brillantor de les estrelles del cel amb la brillantor que indiquen els mapes, gent de tot el món aprendran com els llums de la seva zona contribueixen a augmentar la contaminació lumínica. Les vostres aportacions a la base de dades activa faran palesa la visibilitat del cel nocturn."

        $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>' + $newText + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

        $rng.InsertXML($xml)
        break
    }
}
